$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F21").Value = "93_referral_statement"
$ws.Range("F35").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F37").Value = "ppe"
$ws.Range("F38").Value = "ppe"
$ws.Range("F42").Value = "application instructions"
$ws.Range("F45").Value = "env warning - water || off target movement"
$ws.Range("F46").Value = "env warning - water"
$ws.Range("F48").Value = "off target movement"
$ws.Range("F55").Value = "application instructions"
$ws.Range("F56").Value = "application instructions"
$ws.Range("F57").Value = "application instructions"
$ws.Range("F58").Value = "application instructions"
$ws.Range("F59").Value = "134_non-agriculture_use_requirements"
$ws.Range("F60").Value = "135_product_information"
$ws.Range("F61").Value = "135_product_information"
$ws.Range("F62").Value = "use restrictions"
$ws.Range("F66").Value = "mixing"
$ws.Range("F67").Value = "mixing"
$ws.Range("F68").Value = "mixing"
$ws.Range("F69").Value = "application instructions || off target movement"
$ws.Range("F70").Value = "application instructions"
$ws.Range("F71").Value = "application instructions"
$ws.Range("F72").Value = "application instructions"
$ws.Range("F73").Value = "use restrictions"
$ws.Range("F76").Value = "application instructions"
$ws.Range("F185").Value = "application instructions"
$ws.Range("F191").Value = "application instructions"
$ws.Range("F193").Value = "154_pesticide_storage"
